# Added parallel execution and updated excel reader
# Append new test scenario rows to the sheet (rows 3-5)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Unsuccessful login with invalid credentials
$ws.Range("A3").Value = "Unsuccessfull login with invalid credentials"
$ws.Range("B3").Value = "testingInvald"
$ws.Range("C3").Value = "cvjdvjjvj"

# Row 4: Login with empty fields (only Scenario Name populated)
$ws.Range("A4").Value = "Login with empty fields"

# Row 5: User log out successfully
$ws.Range("A5").Value = "User log out Successfully"
$ws.Range("B5").Value = "Admin"
$ws.Range("C5").Value = "admin123"

# Update the active selection to reflect the last edited cell
$ws.Range("C5").Select() | Out-Null
